$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text to preserve formatting (e.g. "1.00", "600.72")
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.735.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.678.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.402'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.159.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.586.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.673.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000111'
$ws.Range("D24").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '530.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '158.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '164.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.15'
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) column (E) - preserve leading/trailing spaces via explicit text
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +6.16%  '
$ws.Range("E9").Value = '  +5.15%  '
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("E11").Value = '  -3.76%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("E13").Value = '  -2.51%  '
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("E24").Value = '  +5.17%  '
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("E28").Value = '  -6.02%  '
$ws.Range("E29").Value = '  -2.33%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("E34").Value = '  -3.83%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("E44").Value = '  +2.35%  '
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("E48").Value = '  -3.13%  '
$ws.Range("E49").Value = '  +12.69%  '
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("E51").Value = '  -4.05%  '
